$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("NBR", "BAR")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the first 4 data rows (old Cutoff 0-3), shifting the remaining
    # rows (old Cutoff 4-18) up into rows 2-16.
    $ws.Range("A2:A5").EntireRow.Delete()

    # Renumber the Cutoff column (A) back to a contiguous 0-based sequence
    # for the rows that remain (B and C keep the values they carried up).
    for ($i = 0; $i -le 14; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
